# Fill in the "Measured Value" (F) and "True Value" (G) columns on the
# "Table 4.1" sheet. These were previously blank; the MEDIAN formulas in
# column H (and everything downstream on "Table 4.2" + the charts) pick
# up the new values automatically on recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 4.1")

$ws.Range("F3").Value = 41
$ws.Range("G3").Value = 40.01

$ws.Range("F4").Value = 29
$ws.Range("G4").Value = 27.09

$ws.Range("F5").Value = 16
$ws.Range("G5").Value = 16.04

$ws.Range("F6").Value = 83
$ws.Range("G6").Value = 82.05

$ws.Range("F7").Value = 138.4
$ws.Range("G7").Value = 137.06

$ws.Range("F8").Value = 277.5
$ws.Range("G8").Value = 260.02

$ws.Range("F9").Value = 7.08
$ws.Range("G9").Value = 7.01

$ws.Range("F10").Value = 44.2
$ws.Range("G10").Value = 45

$ws.Range("F11").Value = 105.6
$ws.Range("G11").Value = 103

$ws.Range("F12").Value = 199.7
$ws.Range("G12").Value = 199

$wb.Save()
